# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (column E) and corresponding "Valor Mora" (column F)
# for rows 16-60 are refreshed: the period list now runs in ascending
# order (1607 .. 2003) instead of descending (2003 .. 1607), and the
# "Valor Mora" amount follows the period (29600 for periods <= 1808,
# 31249 for periods >= 1809).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered list of periods (ascending) that now populate column E, rows 16-60.
$periodos = @(
    "1607","1608","1609","1610","1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$firstRow = 16

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $firstRow + $i
    $periodo = $periodos[$i]

    # Column E keeps its existing "Text" number format, but force it
    # explicitly so the period code (e.g. "1607") is not reinterpreted
    # as a number.
    $eCell = $ws.Range("E$row")
    $eCell.NumberFormat = "@"
    $eCell.Value = $periodo

    # Column F: the "mora" amount tracks the period - 29600 through
    # 1808, 31249 from 1809 onward.
    $fCell = $ws.Range("F$row")
    if ([int]$periodo -le 1808) {
        $fCell.Value = 29600
    } else {
        $fCell.Value = 31249
    }
}
